# Sprint 3 Burndown Chart - add a "Comments" column (D) with sprint update
# notes, per "Additional files from the completion of release two".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("D30").Value = "Comments"

# Individual comments, entered in the order they appear in the shared
# string table of the target workbook (header, then story/task comments
# top-to-bottom as authored).
$ws.Range("D38").Value = "Blake: Completed user story 26"
$ws.Range("D33").Value = "Blake: Completed task 1 of story 10"
$ws.Range("D36").Value = "Blake: Completed task 2 of story 10"
$ws.Range("D43").Value = "Ashleigh: Completed task 3 and 4 of additional tasks"
$ws.Range("D37").Value = "Alex: Completed task 1 and 2 of additional tasks"

# Widen the new column so the comments are readable (closest attainable
# width to the authored 33.875 "characters" given this host's fixed
# 6-pt-per-character column metric).
$ws.Columns("D").ColumnWidth = 33

# Update the view's current selection to match the author's last position.
$ws.Range("F41").Select()
